$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the first sheet (Tabla_1), which contained the erroneous
# Author/Name/Type table. The remaining sheets shift into its place.
$wb.Worksheets.Item("Tabla_1").Delete()

# Renumber the remaining sheets back to Tabla_1, Tabla_2, Tabla_3 so the
# sheet names stay sequential after the deletion.
$wb.Worksheets.Item("Tabla_2").Name = "Tabla_1"
$wb.Worksheets.Item("Tabla_3").Name = "Tabla_2"
$wb.Worksheets.Item("Tabla_4").Name = "Tabla_3"
